$wb = $excel.ActiveWorkbook

# --- Update "Logs" worksheet: append row 4 with new test mail data ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A4").Value = "Heb je de CE-certificaten van dit product?"
$logs.Range("B4").Value = "mailmind.test@zohomail.eu"
$logs.Range("C4").Value = "Testmail #14: Heb je de CE-certificaten van dit product?"
$logs.Range("D4").Value = "Productinformatie"
$logs.Range("E4").Value = "Geachte afzender,`nDank u voor uw bericht. Helaas kan ik u op dit moment niet voorzien van de CE-certificaten van het genoemde product. Om u verder te helpen, zou u ons de naam van het specifieke product kunnen doorgeven, zodat we uw vraag nauwkeuriger kunnen beantwoorden.`nIk kijk uit naar uw reactie.`nMet vriendelijke groet,`n[Naam]`nNederlandse e-mailassistent"
$logs.Range("F4").Value = "2025-08-02 00:05:16"
$logs.Range("G4").Value = "Ja"
$logs.Range("H4").Value = "Nee"
$logs.Range("I4").Value = "Ja"
$logs.Range("J4").Value = "Nee"

# --- Extend conditional formatting ranges down to the new row 4 ---
$ranges = @("D2:D3", "G2:G3", "H2:H3", "I2:I3", "J2:J3")
foreach ($addr in $ranges) {
    $col = $addr.Substring(0, 1)
    $newAddr = "$($col)2:$($col)4"
    $fcs = $logs.Range($addr).FormatConditions
    $newRange = $logs.Range($newAddr)
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- Update "Dashboard" worksheet: swap categorie ordering & counts ---
$dashboard = $wb.Worksheets.Item("Dashboard")

$dashboard.Range("A2").Value = "Productinformatie"
$dashboard.Range("B2").Value = 2
$dashboard.Range("A3").Value = "Retour / Terugbetaling"
$dashboard.Range("B3").Value = 1
